# "Generate Report for Handback"
# The 827253f4-...md and 927edbff-...md files have now been handed back
# (both zh-cn and de-de). Update the Overview sheet's per-language status,
# and fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# 827253f4 row (row 4) and the 927edbff row (row 5).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $handedBack
$wsOverview.Range("F4").Value = $handedBack
$wsOverview.Range("E5").Value = $handedBack
$wsOverview.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4 : 827253f4-11c5-4be0-abec-9e127905240a.md
$wsZhCn.Range("C4").Value = $handedBack
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/827253f4-11c5-4be0-abec-9e127905240a.md",
    [Type]::Missing,
    [Type]::Missing,
    "827253f4-11c5-4be0-abec-9e127905240a.md"
)
$wsZhCn.Range("J4").Value = "827253f4-11c5-4be0-abec-9e127905240a.2b4474c7cb471b9fac307285952d49758a6de331.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-13 22:33:28"

# Row 5 : 927edbff-a883-4087-ad65-5f4b84f07fa9.md
$wsZhCn.Range("C5").Value = $handedBack
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I5"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/927edbff-a883-4087-ad65-5f4b84f07fa9.md",
    [Type]::Missing,
    [Type]::Missing,
    "927edbff-a883-4087-ad65-5f4b84f07fa9.md"
)
$wsZhCn.Range("J5").Value = "927edbff-a883-4087-ad65-5f4b84f07fa9.883ddb5e5b94b18e8bf22521251b46e68654aa3c.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-13 22:33:28"

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 : 827253f4-11c5-4be0-abec-9e127905240a.md
$wsDeDe.Range("C4").Value = $handedBack
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/827253f4-11c5-4be0-abec-9e127905240a.md",
    [Type]::Missing,
    [Type]::Missing,
    "827253f4-11c5-4be0-abec-9e127905240a.md"
)
$wsDeDe.Range("J4").Value = "827253f4-11c5-4be0-abec-9e127905240a.2b4474c7cb471b9fac307285952d49758a6de331.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-13 22:33:38"

# Row 5 : 927edbff-a883-4087-ad65-5f4b84f07fa9.md
$wsDeDe.Range("C5").Value = $handedBack
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I5"),
    "https://github.com/OpenLocalizationTestOrg/oltest/blob/116024577fc572e872ccc58bcc881c42fa2fe840/e2e/927edbff-a883-4087-ad65-5f4b84f07fa9.md",
    [Type]::Missing,
    [Type]::Missing,
    "927edbff-a883-4087-ad65-5f4b84f07fa9.md"
)
$wsDeDe.Range("J5").Value = "927edbff-a883-4087-ad65-5f4b84f07fa9.883ddb5e5b94b18e8bf22521251b46e68654aa3c.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-13 22:33:38"
